# Updates cryptos list values (Coin/Link/Price/Volume(1h)) to match the
# refreshed data pulled by the scheduled GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.461.14"
$ws.Range("E2").Value = "  +0.17%  "
# Row 3
$ws.Range("D3").Value = "1.626.31"
$ws.Range("E3").Value = "  +0.67%  "
# Row 4
$ws.Range("D4").Value = "'0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "
# Row 5
$ws.Range("D5").Value = "'213.33"
$ws.Range("D5").Style = "Normal"
# Row 6
$ws.Range("D6").Value = "'0.502"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.12%  "
# Row 7
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.09%  "
# Row 8
$ws.Range("E8").Value = "  +0.08%  "
# Row 9
$ws.Range("D9").Value = "'0.0611"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.21%  "
# Row 10
$ws.Range("D10").Value = "'19.24"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.23%  "
# Row 11
$ws.Range("D11").Value = "'0.0850"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.45%  "
# Row 12
$ws.Range("D12").Value = "1.849.01"
$ws.Range("E12").Value = "  +0.48%  "
# Row 13
$ws.Range("D13").Value = "1.618.55"
$ws.Range("E13").Value = "  -0.63%  "
# Row 14
$ws.Range("E14").Value = "  -0.10%  "
# Row 15
$ws.Range("D15").Value = "'0.513"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.17%  "
# Row 16
$ws.Range("D16").Value = "'63.95"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.15%  "
# Row 17
$ws.Range("D17").Value = "'236.60"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.02%  "
# Row 18
$ws.Range("D18").Value = "26.457.08"
$ws.Range("E18").Value = "  +0.15%  "
# Row 19
$ws.Range("E19").Value = "  +4.20%  "
# Row 20
$ws.Range("D20").Value = "0.0₃0727"
$ws.Range("E20").Value = "  -0.04%  "
# Row 21
$ws.Range("E21").Value = "  +0.10%  "
# Row 22
$ws.Range("D22").Value = "'4.33"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.85%  "
# Row 23
$ws.Range("B23").Value = "Avalanche"
$ws.Range("C23").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D23").Value = "'9.14"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.54%  "
# Row 24
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").Value = "'2.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.46%  "
# Row 25
$ws.Range("D25").Value = "'147.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.49%  "
# Row 26
$ws.Range("E26").Value = "  +0.03%  "
# Row 27
$ws.Range("D27").Value = "'7.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.75%  "
# Row 28
$ws.Range("D28").Value = "'0.114"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.05%  "
# Row 29
$ws.Range("D29").Value = "'15.61"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.87%  "
# Row 30
$ws.Range("D30").Value = "'0.0498"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.05%  "
# Row 31
$ws.Range("E31").Value = "  -0.44%  "
# Row 32
$ws.Range("D32").Value = "1.523.86"
$ws.Range("E32").Value = "  +5.21%  "
# Row 33
$ws.Range("E33").Value = "  +1.22%  "
# Row 34
$ws.Range("E34").Value = "  -0.19%  "
# Row 35
$ws.Range("E35").Value = "  +3.02%  "
# Row 36
$ws.Range("E36").Value = "  +0.01%  "
# Row 37
$ws.Range("D37").Value = "'0.571"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.12%  "
# Row 38
$ws.Range("E38").Value = "  +0.11%  "
# Row 39
$ws.Range("D39").Value = "'0.836"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.03%  "
# Row 40
$ws.Range("E40").Value = "  +0.72%  "
# Row 41
$ws.Range("E41").Value = "  +0.05%  "
# Row 42
$ws.Range("E42").Value = "  +0.50%  "
# Row 43
$ws.Range("D43").Value = "1.761.47"
$ws.Range("E43").Value = "  +0.59%  "
# Row 44
$ws.Range("D44").Value = "'62.95"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.72%  "
# Row 45
$ws.Range("E45").Value = "  +0.23%  "
# Row 46
$ws.Range("D46").Value = "'0.912"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.10%  "
# Row 47
$ws.Range("D47").Value = "'90.63"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.99%  "
# Row 48
$ws.Range("E48").Value = "  +1.32%  "
# Row 49
$ws.Range("E49").Value = "  -0.17%  "
# Row 50
$ws.Range("D50").Value = "'0.0967"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.60%  "
# Row 51
$ws.Range("D51").Value = "'7.49"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.15%  "
